# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) for the affected Leve rows across all job sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 19
$ws.Range("H19").Value = 396.5
$ws.Range("I19").Value = 440.07693
$ws.Range("J19").Value = 358.73334
$ws.Range("K19").Value = 440.07693
$ws.Range("L19").Value = 358.73334
$ws.Range("M19").Value = -265.07693
$ws.Range("N19").Value = -708.73334

# Row 96
$ws.Range("H96").Value = 1351.8
$ws.Range("I96").Value = 1855.3334
$ws.Range("J96").Value = 596.5
$ws.Range("K96").Value = 5566.0002
$ws.Range("L96").Value = 1789.5
$ws.Range("M96").Value = -4193.0002
$ws.Range("N96").Value = -4535.5

# Row 137
$ws.Range("H137").Value = 2160.5098
$ws.Range("I137").Value = 1701.2
$ws.Range("J137").Value = 2602.1538
$ws.Range("K137").Value = 5103.6
$ws.Range("L137").Value = 7806.4614
$ws.Range("M137").Value = -2553.6
$ws.Range("N137").Value = -12906.4614

# Row 138
$ws.Range("H138").Value = 2240.31
$ws.Range("J138").Value = 2404.1477
$ws.Range("L138").Value = 7212.4431
$ws.Range("N138").Value = -17492.4431

# Row 141
$ws.Range("H141").Value = 1094
$ws.Range("I141").Value = 1094
$ws.Range("K141").Value = 3282
$ws.Range("M141").Value = 1898

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 10440.52
$ws.Range("I32").Value = 7688.5757
$ws.Range("K32").Value = 7688.5757
$ws.Range("M32").Value = -7401.5757

# Row 61
$ws.Range("H61").Value = 76924520
$ws.Range("J61").Value = 1587.3334
$ws.Range("L61").Value = 1587.3334
$ws.Range("N61").Value = -2011.3334

# Row 74
$ws.Range("H74").Value = 2360.923
$ws.Range("I74").Value = 1808.3636
$ws.Range("K74").Value = 1808.3636
$ws.Range("M74").Value = -934.3635999999999

# Row 77
$ws.Range("H77").Value = 2360.923
$ws.Range("I77").Value = 1808.3636
$ws.Range("K77").Value = 9041.817999999999
$ws.Range("M77").Value = -4673.817999999999

# Row 97
$ws.Range("H97").Value = 552.25
$ws.Range("I97").Value = 528.7917
$ws.Range("K97").Value = 528.7917
$ws.Range("M97").Value = -32.79169999999999

# Row 136
$ws.Range("H136").Value = 76924520
$ws.Range("J136").Value = 1587.3334
$ws.Range("L136").Value = 4762.0002
$ws.Range("N136").Value = -9862.0002

$ws = $wb.Worksheets.Item("BSM")
# Row 81
$ws.Range("H81").Value = 19730.4
$ws.Range("J81").Value = 19730.4
$ws.Range("L81").Value = 19730.4
$ws.Range("N81").Value = -21852.4

# Row 84
$ws.Range("H84").Value = 19730.4
$ws.Range("J84").Value = 19730.4
$ws.Range("L84").Value = 59191.2
$ws.Range("N84").Value = -69799.20000000001

# Row 96
$ws.Range("H96").Value = 10428
$ws.Range("I96").Value = 10428
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 10428
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -7682
$ws.Range("N96").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 1570.0182
$ws.Range("I31").Value = 1443.0613
$ws.Range("K31").Value = 1443.0613
$ws.Range("M31").Value = -1148.0613

# Row 34
$ws.Range("H34").Value = 1570.0182
$ws.Range("I34").Value = 1443.0613
$ws.Range("K34").Value = 1443.0613
$ws.Range("M34").Value = -1241.0613

# Row 44
$ws.Range("H44").Value = 1999
$ws.Range("I44").Value = 1999
$ws.Range("K44").Value = 1999
$ws.Range("M44").Value = -1557

# Row 97
$ws.Range("H97").Value = 29000
$ws.Range("J97").Value = 29000
$ws.Range("L97").Value = 29000
$ws.Range("N97").Value = -30982

# Row 134
$ws.Range("H134").Value = 16668205
$ws.Range("I134").Value = 1579.8695
$ws.Range("J134").Value = 71429976
$ws.Range("K134").Value = 4739.6085
$ws.Range("L134").Value = 214289928
$ws.Range("M134").Value = -2204.6085
$ws.Range("N134").Value = -214294998

$ws = $wb.Worksheets.Item("CUL")
# Row 39
$ws.Range("H39").Value = 4269.4
$ws.Range("J39").Value = 4386.75
$ws.Range("L39").Value = 13160.25
$ws.Range("N39").Value = -13748.25

# Row 44
$ws.Range("H44").Value = 679.3
$ws.Range("I44").Value = 298.6
$ws.Range("J44").Value = 1060
$ws.Range("K44").Value = 895.8000000000001
$ws.Range("L44").Value = 3180
$ws.Range("M44").Value = -497.8000000000001
$ws.Range("N44").Value = -3976

# Row 46
$ws.Range("H46").Value = 500
$ws.Range("I46").Value = 500
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1500
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -1409
$ws.Range("N46").ClearContents()

# Row 131
$ws.Range("H131").Value = 18897102
$ws.Range("I131").Value = 66667204
$ws.Range("J131").Value = 40482.816
$ws.Range("K131").Value = 200001612
$ws.Range("L131").Value = 121448.448
$ws.Range("M131").Value = -199996572
$ws.Range("N131").Value = -131528.448

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 3071.9487
$ws.Range("I132").Value = 2892.8965
$ws.Range("K132").Value = 8678.6895
$ws.Range("M132").Value = -6148.6895

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 254.3871
$ws.Range("I55").Value = 201.36
$ws.Range("K55").Value = 201.36
$ws.Range("M55").Value = -28.36000000000001

# Row 132
$ws.Range("H132").Value = 2627.3076
$ws.Range("I132").Value = 3009.8
$ws.Range("K132").Value = 9029.400000000001
$ws.Range("M132").Value = -6499.400000000001

# Row 136
$ws.Range("H136").Value = 1673.091
$ws.Range("I136").Value = 942.7143
$ws.Range("J136").Value = 2951.25
$ws.Range("K136").Value = 2828.1429
$ws.Range("L136").Value = 8853.75
$ws.Range("M136").Value = -278.1428999999998
$ws.Range("N136").Value = -13953.75

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 38198
$ws.Range("J46").Value = 38198
$ws.Range("L46").Value = 38198
$ws.Range("N46").Value = -38660

# Row 81
$ws.Range("H81").Value = 1683.6428
$ws.Range("J81").Value = 1854.6364
$ws.Range("L81").Value = 3709.2728
$ws.Range("N81").Value = -5831.272800000001

# Row 84
$ws.Range("H84").Value = 1683.6428
$ws.Range("J84").Value = 1854.6364
$ws.Range("L84").Value = 18546.364
$ws.Range("N84").Value = -29154.364

# Row 121
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# Row 134
$ws.Range("H134").Value = 38198
$ws.Range("J134").Value = 38198
$ws.Range("L134").Value = 114594
$ws.Range("N134").Value = -119664

# Row 136
$ws.Range("H136").Value = 1244.4783
$ws.Range("I136").Value = 919.05554
$ws.Range("J136").Value = 2416
$ws.Range("K136").Value = 2757.16662
$ws.Range("L136").Value = 7248
$ws.Range("M136").Value = -207.16662
$ws.Range("N136").Value = -12348

# Row 141
$ws.Range("H141").Value = 49052.777
$ws.Range("J141").Value = 49052.777
$ws.Range("L141").Value = 49052.777
$ws.Range("N141").Value = -59412.777
